# Fix typo: "asterix" -> "asterisk", and update the "wrong" icon label to "_test_"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("E5").Value = "asterisk"
$ws.Range("E14").Value = "_test_"
